# Update WanHaoBillDetails: add entries dated 2017-10-08 (serial 43016)
# Row 14: income (收入) 生活费 3000
# Row 15: expense (支出) 生活费(10/9~10/23) 500

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14 - income of 3000 for living expenses
$ws.Range("B14").Value = "收入"
$ws.Range("C14").Value = 3000
$ws.Range("D14").Value = 43016
$ws.Range("E14").Value = "生活费"

# Row 15 - expense of 500 for living expenses (10/9~10/23)
$ws.Range("B15").Value = "支出"
$ws.Range("C15").Value = 500
$ws.Range("D15").Value = 43016
$ws.Range("E15").Value = "生活费(10/9~10/23)"

# Copy the date style (numFmtId 14, same as other rows) onto the new date
# cells instead of re-declaring the number format, which would otherwise
# register a brand-new custom numFmt entry.
$ws.Range("D13").Copy()
$ws.Range("D14:D15").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Match the author's final selection
$ws.Range("J13").Select()
